$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest values.
# D-column values are forced to Text format so numeric-looking strings (e.g. "20.70",
# "0.06720", "287.62") keep their exact original textual representation instead of
# being coerced into numbers by Excel (which would strip trailing zeros, etc).
# The cell Style is reset to "Normal" afterwards so no stray formatting/style index
# is left behind on cells that originally had none.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.427.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.570.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3692"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.42"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.17%  "
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07492"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.926"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.885"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.559.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001112"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06720"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.419"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.414.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.616"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.931"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.736.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.083"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.972"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02427"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06370"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.297"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2204"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.316"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6026"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.775"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.033"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.193"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07188"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
